$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21: bold header labels (mirrors the "title" style used by the other
# section headers at B11:D11 / B17:D17)
$ws.Range("B21").Value = "Number of employees"
$ws.Range("C21").Value = "Assets (local currency, unless noted otherwise)"
$ws.Range("D21").Value = "Turnover (local currency, unless noted otherwise)"
$ws.Range("B21:D21").Font.Bold = $true

# Rows 22-25: MSME size categories (Micro/Small/Medium/Large) with blank
# data columns B:D ready for figures to be filled in later.
$ws.Range("A22").Value = "Micro"
$ws.Range("B22").Value = ""
$ws.Range("C22").Value = ""
$ws.Range("D22").Value = ""

$ws.Range("A23").Value = "Small"
$ws.Range("B23").Value = ""
$ws.Range("C23").Value = ""
$ws.Range("D23").Value = ""

$ws.Range("A24").Value = "Medium"
$ws.Range("B24").Value = ""
$ws.Range("C24").Value = ""
$ws.Range("D24").Value = ""

$ws.Range("A25").Value = "Large"
$ws.Range("B25").Value = ""
$ws.Range("C25").Value = ""
$ws.Range("D25").Value = ""
